# Update column F ("dSF") values for the rows that changed on repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -5
    8  = -2
    9  = -2
    11 = 5
    12 = -4
    13 = -3
    15 = 0
    22 = -2
    27 = -6
    28 = -4
    29 = -6
    35 = -4
    40 = -3
    42 = -2
    46 = -7
    49 = -4
    50 = -3
    51 = -6
    52 = -6
    53 = -3
    54 = -1
    56 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
